$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Part names for the new rows, entered first (matches author's shared-string order)
$ws.Range("A7").Value = "IMU"
$ws.Range("A8").Value = "STM32 Nucleo Board"
$ws.Range("A9").Value = "Caster Wheel"
$ws.Range("A10").Value = "OLED Screen"

# Row 7 - IMU
$ws.Range("B7").Value = 1
$ws.Range("D7").Value = 19.95
$ws.Range("D7").NumberFormat = '"$"#,##0.00'
$ws.Range("H7").Formula = "=(B7+C7)*D7"
$ws.Range("H7").NumberFormat = '"$"#,##0.00'
$ws.Range("I7").Value = "https://www.adafruit.com/product/4502"

# Row 8 - STM32 Nucleo Board
$ws.Range("B8").Value = 1

# Row 9 - Caster Wheel (no other details)

# Row 10 - OLED Screen
$ws.Range("B10").Value = 1
$ws.Range("D10").Value = 17.5
$ws.Range("I10").Value = "https://www.adafruit.com/product/661"

# Widen column A to fit the new longer part names
$ws.Columns.Item(1).ColumnWidth = 16.67

# Match the cursor/selection position left by the author
$ws.Range("G9").Select()
